# Applies the F-column count increments and the row 31-36 content shift
# described by the diff, to both '展览' (sheet 1) and '全部类型' (sheet 4).

$wb = $excel.ActiveWorkbook

$simpleF = @{
    6 = 125
    7 = 1230
    8 = 1519
    9 = 334
    10 = 378
    12 = 142
    15 = 104
    16 = 270
    19 = 1714
    23 = 658
    26 = 4127
    28 = 477
    30 = 1075
}

$targetSheets = @(1, 4)

foreach ($sheetIndex in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    foreach ($row in $simpleF.Keys) {
        $ws.Cells.Item($row, 6).Value = $simpleF[$row]
    }

    # Row 31
    $ws.Cells.Item(31, 3).Value = "景德镇·原神X崩铁X崩坏动漫展only"
    $ws.Cells.Item(31, 4).Value = "陶阳南路188号 晨枫臻品酒店"
    $ws.Cells.Item(31, 5).Value = "2024.03.16 10:00-03.16 17:00"
    $ws.Cells.Item(31, 6).Value = 46
    $ws.Cells.Item(31, 7).Value = 55
    $ws.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80920"
    $ws.Cells.Item(31, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"

    # Row 32
    $ws.Cells.Item(32, 3).Value = "江西·ShiningStaR动漫游戏文化节5th"
    $ws.Cells.Item(32, 4).Value = "高新开发区紫阳大道666号 江西奥林匹克体育中心综合训练馆"
    $ws.Cells.Item(32, 5).Value = "2024.03.16 09:30-03.17 17:00"
    $ws.Cells.Item(32, 6).Value = 499
    $ws.Cells.Item(32, 7).Value = 60
    $ws.Cells.Item(32, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81792"
    $ws.Cells.Item(32, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/2l16aHBJ1707209383729.jpeg"

    # Row 33
    $ws.Cells.Item(33, 2).NumberFormat = "@"
    $ws.Cells.Item(33, 2).Value = "2024-03-23"
    $ws.Cells.Item(33, 3).Value = "上饶·原×铁×崩only"
    $ws.Cells.Item(33, 4).Value = "五三东大道42号 回禾酒店"
    $ws.Cells.Item(33, 5).Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Cells.Item(33, 6).Value = 22
    $ws.Cells.Item(33, 7).Value = 60
    $ws.Cells.Item(33, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81103"
    $ws.Cells.Item(33, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"

    # Row 34
    $ws.Cells.Item(34, 3).Value = "南昌·AP动漫游戏嘉年华"
    $ws.Cells.Item(34, 4).Value = "八一桥街道青山南路118号 蓝海会展中心"
    $ws.Cells.Item(34, 5).Value = "2024.03.23 09:00-03.24 17:00"
    $ws.Cells.Item(34, 6).Value = 226
    $ws.Cells.Item(34, 7).Value = 60
    $ws.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81232"
    $ws.Cells.Item(34, 9).Value = "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"

    # Row 35
    $ws.Cells.Item(35, 3).Value = "南昌·原X穹X崩only"
    $ws.Cells.Item(35, 4).Value = "丰和北大道299号 新吉花园酒店"
    $ws.Cells.Item(35, 5).Value = "2024.03.23 10:00-03.23 17:00"
    $ws.Cells.Item(35, 6).Value = 48
    $ws.Cells.Item(35, 7).Value = 65
    $ws.Cells.Item(35, 8).Value = "https://show.bilibili.com/platform/detail.html?id=80807"
    $ws.Cells.Item(35, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/rY4v2Opx1705051458246.jpeg"

    # Row 36
    $ws.Cells.Item(36, 2).NumberFormat = "@"
    $ws.Cells.Item(36, 2).Value = "2024-03-30"
    $ws.Cells.Item(36, 3).Value = "南昌·CM01动漫游戏博览会"
    $ws.Cells.Item(36, 4).Value = "怀玉山大道1315号 南昌绿地国际博览中心"
    $ws.Cells.Item(36, 5).Value = "2024.03.30 10:00-03.31 17:00"
    $ws.Cells.Item(36, 6).Value = 133
    $ws.Cells.Item(36, 7).Value = 55
    $ws.Cells.Item(36, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81691"
    $ws.Cells.Item(36, 9).Value = "//i2.hdslb.com/bfs/openplatform/202402/IYLaH7AS1706866218597.png"

}

Write-Host "Done applying edits"